$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.492.72"
$ws.Range("E2").Value = "'  -0.53%  "

$ws.Range("D3").Value = "'1.834.53"

$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = "'  +0.54%  "

$ws.Range("D5").Value = "'313.41"
$ws.Range("E5").Value = "'  +0.03%  "

$ws.Range("D6").Value = "'1.005"
$ws.Range("E6").Value = "'  +0.39%  "

$ws.Range("D7").Value = "'0.4236"
$ws.Range("E7").Value = "'  +0.05%  "

$ws.Range("D8").Value = "'0.3657"
$ws.Range("E8").Value = "'  +0.44%  "

$ws.Range("D9").Value = "'0.07217"
$ws.Range("E9").Value = "'  -1.06%  "

$ws.Range("D10").Value = "'0.8648"
$ws.Range("E10").Value = "'  -1.26%  "

$ws.Range("D11").Value = "'20.72"
$ws.Range("E11").Value = "'  +0.16%  "

$ws.Range("D12").Value = "'1.876.28"
$ws.Range("E12").Value = "'  +1.28%  "

$ws.Range("D13").Value = "'5.374"
$ws.Range("E13").Value = "'  +0.64%  "

$ws.Range("D14").Value = "'6.499"
$ws.Range("E14").Value = "'  -0.39%  "

$ws.Range("D15").Value = "'0.06965"
$ws.Range("E15").Value = "'  +1.40%  "

$ws.Range("D16").Value = "'1.007"
$ws.Range("E16").Value = "'  +0.43%  "

$ws.Range("D17").Value = "'79.52"
$ws.Range("E17").Value = "'  -0.12%  "

$ws.Range("D18").Value = "'0.000008982"
$ws.Range("E18").Value = "'  +0.97%  "

$ws.Range("D19").Value = "'1.008"
$ws.Range("E19").Value = "'  +0.73%  "

$ws.Range("D20").Value = "'15.39"
$ws.Range("E20").Value = "'  +0.22%  "

$ws.Range("D21").Value = "'27.895.53"
$ws.Range("E21").Value = "'  +0.81%  "

$ws.Range("D22").Value = "'5.031"
$ws.Range("E22").Value = "'  +0.83%  "

$ws.Range("D23").Value = "'10.77"
$ws.Range("E23").Value = "'  +3.77%  "

$ws.Range("D24").Value = "'2.138.01"
$ws.Range("E24").Value = "'  +1.81%  "

$ws.Range("D25").Value = "'1.964"
$ws.Range("E25").Value = "'  -0.93%  "

$ws.Range("D26").Value = "'154.10"
$ws.Range("E26").Value = "'  +0.47%  "

$ws.Range("E27").Value = "'  -2.76%  "

$ws.Range("D28").Value = "'5.243"
$ws.Range("E28").Value = "'  -0.21%  "

$ws.Range("D29").Value = "'114.66"
$ws.Range("E29").Value = "'  -6.01%  "

$ws.Range("D30").Value = "'1.819"
$ws.Range("E30").Value = "'  -3.11%  "

$ws.Range("D31").Value = "'0.08869"
$ws.Range("E31").Value = "'  +0.06%  "

$ws.Range("D32").Value = "'0.7716"
$ws.Range("E32").Value = "'  +0.54%  "

$ws.Range("D33").Value = "'4.539"
$ws.Range("E33").Value = "'  -0.21%  "

$ws.Range("D34").Value = "'2.955"
$ws.Range("E34").Value = "'  -0.66%  "

$ws.Range("D35").Value = "'1.152"
$ws.Range("E35").Value = "'  +4.08%  "

$ws.Range("D36").Value = "'1.005"
$ws.Range("E36").Value = "'  +0.49%  "

$ws.Range("D37").Value = "'1.097"
$ws.Range("E37").Value = "'  +0.12%  "

$ws.Range("D38").Value = "'0.05367"
$ws.Range("E38").Value = "'  +0.40%  "

$ws.Range("D39").Value = "'0.01943"
$ws.Range("E39").Value = "'  +0.71%  "

$ws.Range("D40").Value = "'2.826"
$ws.Range("E40").Value = "'  +0.45%  "

$ws.Range("D41").Value = "'0.5116"
$ws.Range("E41").Value = "'  +0.46%  "

$ws.Range("B42").Value = "'FraxShare"
$ws.Range("C42").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'6.814"
$ws.Range("E42").Value = "'  -1.04%  "

$ws.Range("B43").Value = "'Algorand"
$ws.Range("C43").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.1647"
$ws.Range("E43").Value = "'  +0.01%  "

$ws.Range("D44").Value = "'8.461"
$ws.Range("E44").Value = "'  +1.84%  "

$ws.Range("D45").Value = "'10.44"
$ws.Range("E45").Value = "'  +1.22%  "

$ws.Range("D46").Value = "'0.06522"
$ws.Range("E46").Value = "'  -0.25%  "

$ws.Range("B47").Value = "'Quant"
$ws.Range("C47").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "'105.30"
$ws.Range("E47").Value = "'  -0.25%  "

$ws.Range("B48").Value = "'Decentraland"
$ws.Range("C48").Value = "'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "'0.4686"
$ws.Range("E48").Value = "'  +0.07%  "

$ws.Range("D49").Value = "'1.004"
$ws.Range("E49").Value = "'  +0.37%  "

$ws.Range("D50").Value = "'1.621"
$ws.Range("E50").Value = "'  -0.21%  "

$ws.Range("D51").Value = "'1.805"
$ws.Range("E51").Value = "'  +5.19%  "
